# Updates cryptos list values to match the latest scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '27.268.85'
$ws.Range('E2').Value = '  -4.41%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '1.854.88'
$ws.Range('E3').Value = '  -5.60%  '

# Row 4: TetherUSD
$ws.Range('E4').Value = '  -1.07%  '

# Row 5: BNB
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '321.01'
$ws.Range('E5').Value = '  -0.68%  '

# Row 6: USDC
$ws.Range('E6').Value = '  -1.04%  '

# Row 7: XRP
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4489'
$ws.Range('E7').Value = '  -5.76%  '

# Row 8: Cardano
$ws.Range('E8').Value = '  -5.06%  '

# Row 9: OKB
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '47.66'
$ws.Range('E9').Value = '  -11.71%  '

# Row 10: Dogecoin
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.07883'
$ws.Range('E10').Value = '  -7.02%  '

# Row 11: Polygon
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '1.017'
$ws.Range('E11').Value = '  -4.08%  '

# Row 12: Solana
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '21.32'
$ws.Range('E12').Value = '  -4.78%  '

# Row 13: WrappedEther
$ws.Range('D13').Value = '1.849.23'
$ws.Range('E13').Value = '  -6.86%  '

# Row 14: Chainlink
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '7.152'
$ws.Range('E14').Value = '  -5.91%  '

# Row 15: Polkadot
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '5.864'
$ws.Range('E15').Value = '  -5.25%  '

# Row 16: BinanceUSD
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '1.002'
$ws.Range('E16').Value = '  -1.16%  '

# Row 17: ShibaInu
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.00001029'
$ws.Range('E17').Value = '  -4.05%  '

# Row 18: Litecoin
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '85.56'
$ws.Range('E18').Value = '  -5.90%  '

# Row 19: TRON
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06525'
$ws.Range('E19').Value = '  -1.52%  '

# Row 20: Avalanche
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '16.90'
$ws.Range('E20').Value = '  -8.77%  '

# Row 21: Dai
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.9997'
$ws.Range('E21').Value = '  -1.20%  '

# Row 22: Uniswap
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.484'
$ws.Range('E22').Value = '  -6.42%  '

# Row 23: WrappedBTC
$ws.Range('D23').Value = '27.265.57'

# Row 24: Cosmos
$ws.Range('E24').Value = '  -6.12%  '

# Row 25: Toncoin
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.263'
$ws.Range('E25').Value = '  -1.66%  '

# Row 26: WrappedliquidstakedEther2.0
$ws.Range('D26').Value = '2.079.52'
$ws.Range('E26').Value = '  -6.45%  '

# Row 27: Monero
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '151.55'
$ws.Range('E27').Value = '  -2.69%  '

# Row 28: EthereumClassic
$ws.Range('E28').Value = '  -3.41%  '

# Row 29: LidoDAOToken
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.058'
$ws.Range('E29').Value = '  -5.04%  '

# Row 30: InternetComputer(DFINITY)
$ws.Range('E30').Value = '  -7.46%  '

# Row 31: BitcoinCash
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '120.20'
$ws.Range('E31').Value = '  -3.57%  '

# Row 32: ImmutableX
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.9366'
$ws.Range('E32').Value = '  -4.67%  '

# Row 33: Stellar
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.09264'
$ws.Range('E33').Value = '  -3.92%  '

# Row 34: ARBITRUM
$ws.Range('E34').Value = '  +0.11%  '

# Row 35: HuobiToken
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '3.567'

# Row 36: Filecoin
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.287'
$ws.Range('E36').Value = '  -6.01%  '

# Row 37: VeChain
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.02221'
$ws.Range('E37').Value = '  -4.84%  '

# Row 38: Hedera
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.05979'
$ws.Range('E38').Value = '  -4.09%  '

# Row 39: TrustWalletToken
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.202'
$ws.Range('E39').Value = '  -4.35%  '

# Row 40: FraxShare
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '8.281'
$ws.Range('E40').Value = '  -9.48%  '

# Row 41: Frax
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.000'
$ws.Range('E41').Value = '  -1.09%  '

# Row 42: TheSandbox
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.5898'
$ws.Range('E42').Value = '  -5.05%  '

# Row 43: Algorand
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.1882'
$ws.Range('E43').Value = '  -1.85%  '

# Row 44: Aptos
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '10.09'
$ws.Range('E44').Value = '  -9.70%  '

# Row 45: WEMIXTOKEN
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.261'
$ws.Range('E45').Value = '  -6.63%  '

# Row 46: Decentraland
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.5613'
$ws.Range('E46').Value = '  -5.64%  '

# Row 47: EnergySwap
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '11.85'
$ws.Range('E47').Value = '  -9.49%  '

# Row 48: NEARProtocol -> PancakeSwap (rows 48/49 swapped rank order)
$ws.Range('B48').Value = 'PancakeSwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '3.350'
$ws.Range('E48').Value = '  -1.85%  '

# Row 49: PancakeSwap -> NEARProtocol (rows 48/49 swapped rank order)
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.916'
$ws.Range('E49').Value = '  -6.92%  '

# Row 50: Cronos
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.06803'
$ws.Range('E50').Value = '  -0.16%  '

# Row 51: Quant
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '108.24'
$ws.Range('E51').Value = '  -2.71%  '
